# XBlastMod register map - apply commit edits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: F3 unchanged text "iSW1/A19" (rewritten so the shared-string edit below
# doesn't change this cell's displayed text)
$ws.Range("F3").Value = "iSW1/A19"

# Row 4: H4 "D0_control" (kept, now shares the string added via H10 below)
$ws.Range("H4").Value = "D0_control"

# Row 10 (0xF70E OS Bank ctrl): new G10/H10 values
$ws.Range("G10").Value = "A19"
$ws.Range("H10").Value = "D0_control"

# The "iSW1/A19" label (bit1 elsewhere) is shortened to "iSW1" -- J3 (bit0 row3)
# and I10 (bit1 row10) both showed that label and now read "iSW1"
$ws.Range("J3").Value = "iSW1"
$ws.Range("I10").Value = "iSW1"

# Remove stray empty cell K21 (row 21 had no other content)
$ws.Range("K21").ClearContents()

# Column I (bit1 column) widened slightly after the content changes above
$ws.Columns.Item(9).ColumnWidth = 20.3

# Restore selection position seen in the saved file
$ws.Range("H19").Select()
